# Rebuild the "UniformF-HW45" sheet with the new simulation-scheme rows
# (HW scheme list) and the reordered [h,k,l] / pairing column headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clear everything beyond column T (old sheet went out to AN) ------
$ws.Range("U1:AN19").Clear()

# --- 2. New column headers for row 2 (C2:T2) ------------------------------
$colHeaders = @(
    "[3, 2, 1]",
    "[3, 1, 0]",
    "[2, 2, 2]",
    "[1, 1, 0]",
    "[2, 0, 0]",
    "[2, 2, 0]",
    "[4, 0, 0]",
    "[2, 1, 1]",
    "1Pair-A",
    "1Pair-B",
    "2Pairs-A",
    "2Pairs-B",
    "3Pairs-A",
    "3Pairs-B",
    "3Pairs-C",
    "4Pairs",
    "5A4F",
    "MaxUnique"
)

for ($i = 0; $i -lt $colHeaders.Length; $i++) {
    $col = 3 + $i   # C = 3
    $ws.Cells.Item(2, $col).Value = $colHeaders[$i]
}

# --- 3. New row labels (column B) for rows 3-29, replacing/extending the --
# --- previous scheme list --------------------------------------------------
$rowLabels = @(
    "Spiral5",
    "RotRing OmegaMax-90",
    "Equal Angle",
    "Tilt Rotate",
    "CLR",
    "Rizzie Hex",
    "Thomas Hex",
    "Tilt Rotate_Partial",
    "RotRing OmegaMax-60",
    "Equal Angle_Partial",
    "Rizzie Hex_Partial",
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris Single",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD",
    "Ring Perpendicular to TD",
    "OffsetFTD",
    "OffsetATD",
    "OffsetF45",
    "OffsetA45",
    "OffsetFRD",
    "OffsetARD",
    "Gaussian Quadrature",
    "Michael-CCHex",
    "Michael-SNHex"
)

for ($i = 0; $i -lt $rowLabels.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 1).Value = $i + 1          # column A index (1..27)
    $ws.Cells.Item($row, 2).Value = $rowLabels[$i]  # column B scheme name

    for ($col = 3; $col -le 20; $col++) {           # C..T = 1
        $ws.Cells.Item($row, $col).Value = 1
    }
}

# --- 4. Apply the row-number style (bold/border/center, same as A2:A19) ---
# --- to the newly added rows 20-29 -----------------------------------------
$ws.Range("A19").Copy()
$ws.Range("A20:A29").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
